$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "... est représentée par Monsieur **El Hadji Mamadou FAYE**, son
#     Directeur Général, "
#  ->
#    "... est représentée par Madame **Jenny MVOU,** son Directeur Général
#     Adjointe, "
# ---------------------------------------------------------------------------

# 1a) Drop "Monsieur " (leaves " est représentée par ").
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Monsieur ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 1b) Swap the bold name.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("El Hadji Mamadou FAYE", $false, $false, $false, $false, $false, $true, 1, $false, "Jenny MVOU,", 2) | Out-Null

# 1c) Insert the non-bold "Madame " right before the (now bold) name, and
#     make sure only the inserted text stays non-bold.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Jenny MVOU,"
$find.Execute() | Out-Null
$startPos = $find.Parent.Start
$insertRange = $d.Range($startPos, $startPos)
$insertRange.InsertBefore("Madame ")
$newRange = $d.Range($startPos, $startPos + 7)
$newRange.Font.Bold = 0

# 1d) Tidy the trailing title.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(", son Directeur Général, ", $false, $false, $false, $false, $false, $true, 1, $false, " son Directeur Général Adjointe, ", 2) | Out-Null

Write-Host "Edit complete."
